$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each data cell in this sheet was written as a plain (non-numeric-typed)
# string by the upstream scraper, even though many values look numeric
# (prices, percentages, the hour counter). Excel.Range.Value auto-detects
# numbers/percentages from plain strings, so we force "Text" number format
# before assigning, then reset the style back to "Normal" so we do not leave
# a stray number-format style attached to the cell (matches original look).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "298.86"
Set-TextValue $ws.Range("E2") "0.86%"
Set-TextValue $ws.Range("G2") "11"
Set-TextValue $ws.Range("D3") "31.36"
Set-TextValue $ws.Range("E3") "0.47%"
Set-TextValue $ws.Range("G3") "11"
Set-TextValue $ws.Range("D4") "5.142"
Set-TextValue $ws.Range("E4") "1.02%"
Set-TextValue $ws.Range("G4") "11"
Set-TextValue $ws.Range("D5") "0.08027"
Set-TextValue $ws.Range("E5") "9.55%"
Set-TextValue $ws.Range("G5") "11"
Set-TextValue $ws.Range("D6") "2.661"
Set-TextValue $ws.Range("E6") "60.74%"
Set-TextValue $ws.Range("G6") "11"
Set-TextValue $ws.Range("D7") "7.846"
Set-TextValue $ws.Range("E7") "1.75%"
Set-TextValue $ws.Range("G7") "11"
Set-TextValue $ws.Range("D8") "3.825"
Set-TextValue $ws.Range("E8") "2.37%"
Set-TextValue $ws.Range("G8") "11"
Set-TextValue $ws.Range("D9") "0.9182"
Set-TextValue $ws.Range("E9") "0.07%"
Set-TextValue $ws.Range("G9") "11"
Set-TextValue $ws.Range("D10") "0.1737"
Set-TextValue $ws.Range("E10") "3.92%"
Set-TextValue $ws.Range("G10") "11"
Set-TextValue $ws.Range("D11") "0.07206"
Set-TextValue $ws.Range("E11") "1.69%"
Set-TextValue $ws.Range("G11") "11"
Set-TextValue $ws.Range("D12") "0.08271"
Set-TextValue $ws.Range("E12") "3.58%"
Set-TextValue $ws.Range("G12") "11"
Set-TextValue $ws.Range("D13") "0.02997"
Set-TextValue $ws.Range("E13") "0.39%"
Set-TextValue $ws.Range("G13") "11"
Set-TextValue $ws.Range("D14") "0.09968"
Set-TextValue $ws.Range("E14") "0.43%"
Set-TextValue $ws.Range("G14") "11"
Set-TextValue $ws.Range("D15") "0.001505"
Set-TextValue $ws.Range("E15") "0.65%"
Set-TextValue $ws.Range("G15") "11"
Set-TextValue $ws.Range("D16") "0.006101"
Set-TextValue $ws.Range("E16") "-0.40%"
Set-TextValue $ws.Range("G16") "11"
Set-TextValue $ws.Range("D17") "3.499"
Set-TextValue $ws.Range("E17") "1.41%"
Set-TextValue $ws.Range("G17") "11"
Set-TextValue $ws.Range("D18") "2.250"
Set-TextValue $ws.Range("E18") "0.98%"
Set-TextValue $ws.Range("G18") "11"
Set-TextValue $ws.Range("E19") "0.35%"
Set-TextValue $ws.Range("G19") "11"
Set-TextValue $ws.Range("D20") "0.1320"
Set-TextValue $ws.Range("E20") "-0.98%"
Set-TextValue $ws.Range("G20") "11"
Set-TextValue $ws.Range("D21") "4.640"
Set-TextValue $ws.Range("E21") "1.83%"
Set-TextValue $ws.Range("G21") "11"
Set-TextValue $ws.Range("D22") "0.1600"
Set-TextValue $ws.Range("E22") "3.23%"
Set-TextValue $ws.Range("G22") "11"
Set-TextValue $ws.Range("D23") "0.04582"
Set-TextValue $ws.Range("E23") "-0.81%"
Set-TextValue $ws.Range("G23") "11"
Set-TextValue $ws.Range("D24") "0.001261"
Set-TextValue $ws.Range("E24") "3.66%"
Set-TextValue $ws.Range("G24") "11"
Set-TextValue $ws.Range("D25") "0.004452"
Set-TextValue $ws.Range("E25") "0.67%"
Set-TextValue $ws.Range("G25") "11"
Set-TextValue $ws.Range("E26") "-9.10%"
Set-TextValue $ws.Range("G26") "11"
Set-TextValue $ws.Range("E27") "83.16%"
Set-TextValue $ws.Range("G27") "11"
Set-TextValue $ws.Range("G28") "11"
Set-TextValue $ws.Range("G29") "11"
Set-TextValue $ws.Range("G30") "11"
Set-TextValue $ws.Range("G31") "11"
Set-TextValue $ws.Range("G32") "11"
Set-TextValue $ws.Range("G33") "11"
Set-TextValue $ws.Range("G34") "11"
Set-TextValue $ws.Range("G35") "11"
Set-TextValue $ws.Range("G36") "11"
Set-TextValue $ws.Range("G37") "11"
Set-TextValue $ws.Range("G38") "11"
Set-TextValue $ws.Range("E39") "9.51%"
Set-TextValue $ws.Range("G39") "11"
Set-TextValue $ws.Range("D40") "0.04514"
Set-TextValue $ws.Range("E40") "2.20%"
Set-TextValue $ws.Range("G40") "11"
Set-TextValue $ws.Range("D41") "0.007034"
Set-TextValue $ws.Range("E41") "-2.52%"
Set-TextValue $ws.Range("G41") "11"
Set-TextValue $ws.Range("D42") "0.1346"
Set-TextValue $ws.Range("E42") "1.39%"
Set-TextValue $ws.Range("G42") "11"
Set-TextValue $ws.Range("D43") "0.002241"
Set-TextValue $ws.Range("E43") "4.81%"
Set-TextValue $ws.Range("G43") "11"
Set-TextValue $ws.Range("D44") "0.01045"
Set-TextValue $ws.Range("E44") "-5.06%"
Set-TextValue $ws.Range("G44") "11"
Set-TextValue $ws.Range("D45") "0.00006483"
Set-TextValue $ws.Range("E45") "7.66%"
Set-TextValue $ws.Range("G45") "11"
Set-TextValue $ws.Range("D46") "0.00000000751"
Set-TextValue $ws.Range("E46") "0.08%"
Set-TextValue $ws.Range("G46") "11"
Set-TextValue $ws.Range("E47") "-39.24%"
Set-TextValue $ws.Range("G47") "11"
Set-TextValue $ws.Range("E48") "15.27%"
Set-TextValue $ws.Range("G48") "11"
Set-TextValue $ws.Range("D49") "0.00002102"
Set-TextValue $ws.Range("E49") "0.08%"
Set-TextValue $ws.Range("G49") "11"
Set-TextValue $ws.Range("D50") "0.0002002"
Set-TextValue $ws.Range("E50") "0.15%"
Set-TextValue $ws.Range("G50") "11"
Set-TextValue $ws.Range("G51") "11"
